$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats / fonts / alignment) from the last
# existing data row down onto the new row so the appended cells reuse the
# workbook's existing styles instead of minting new ones.
$ws.Range("A26:L26").Copy()
$ws.Range("A27:L27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Append a new row of exercise-tracking data at row 27
$ws.Range("A27").Value = 44012
$ws.Range("B27").Value = 98.5
$ws.Range("C27").Value = 106
$ws.Range("D27").Value = 0.93
$ws.Range("E27").Value = "IMPROVED"
$ws.Range("F27").Value = 0.02
$ws.Range("G27").Value = 82.5
$ws.Range("H27").Value = 0.58
$ws.Range("I27").Value = "SAME"
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 28.9
$ws.Range("L27").Value = "OVERWEIGHT"

# Match the saved selection state (whole-sheet selection, no specific active cell)
$ws.Cells.Select()
